# Update the Excel workbook to match the target edit:
# 1) Row with Caso=3441 ("DON PEDRO DE MENDOZA AV. 1487") is removed entirely.
# 2) Row with Caso=5716 ("NECOCHEA 1315") is removed entirely.
# 3) H19's observation text is updated to mention the correction.
# Removing the whole rows shifts all subsequent rows up, which matches the
# diff (every row below the deleted ones moves up by one, or two once both
# deletions are accounted for).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update H19 first (row numbers above the deleted rows are unaffected).
$ws.Range("H19").Value = "Picada cambiaron la incorrecta"

# Delete the lower row first so the earlier row's index doesn't shift before
# we get to it.
$ws.Rows(30).Delete()
$ws.Rows(21).Delete()
